# Auto-generated cell update script
# Applies scheduled-runner price/profit refresh values per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1052.1538
$ws.Range("I19").Value = 967.8
$ws.Range("K19").Value = 967.8
$ws.Range("M19").Value = -792.8
$ws.Range("H61").Value = 556.25
$ws.Range("I61").Value = 556.25
$ws.Range("K61").Value = 1668.75
$ws.Range("M61").Value = -1496.75
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2001
$ws.Range("K62").Value = 2001
$ws.Range("M62").Value = -1377
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2001
$ws.Range("K65").Value = 10005
$ws.Range("M65").Value = -6885
$ws.Range("H100").Value = 2742.8
$ws.Range("I100").Value = 2236.3333
$ws.Range("J100").Value = 3502.5
$ws.Range("K100").Value = 2236.3333
$ws.Range("L100").Value = 3502.5
$ws.Range("M100").Value = -1695.3333
$ws.Range("N100").Value = -4584.5
$ws.Range("H138").Value = 4618.3657
$ws.Range("I138").Value = 1124.375
$ws.Range("J138").Value = 6854.52
$ws.Range("K138").Value = 3373.125
$ws.Range("L138").Value = 20563.56
$ws.Range("M138").Value = 1766.875
$ws.Range("N138").Value = -30843.56

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 126499.836
$ws.Range("I34").Value = 84500
$ws.Range("J34").Value = 147499.75
$ws.Range("K34").Value = 84500
$ws.Range("L34").Value = 147499.75
$ws.Range("M34").Value = -84229
$ws.Range("N34").Value = -148041.75
$ws.Range("H97").Value = 4848.3076
$ws.Range("I97").Value = 4866.1816
$ws.Range("J97").Value = 4750
$ws.Range("K97").Value = 4866.1816
$ws.Range("L97").Value = 4750
$ws.Range("M97").Value = -4370.1816
$ws.Range("N97").Value = -5742
$ws.Range("H102").Value = 6749.5
$ws.Range("I102").Value = 4499
$ws.Range("K102").Value = 4499
$ws.Range("M102").Value = -2877
$ws.Range("H129").Value = 74999.336
$ws.Range("J129").Value = 74999.336
$ws.Range("L129").Value = 74999.336
$ws.Range("N129").Value = -84999.336
$ws.Range("H132").Value = 53293.793
$ws.Range("I132").Value = 6101.45
$ws.Range("J132").Value = 158165.67
$ws.Range("K132").Value = 18304.35
$ws.Range("L132").Value = 474497.01
$ws.Range("M132").Value = -15774.35
$ws.Range("N132").Value = -479557.01

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 38021.637
$ws.Range("I86").Value = 25723.625
$ws.Range("J86").Value = 70816.336
$ws.Range("K86").Value = 25723.625
$ws.Range("L86").Value = 70816.336
$ws.Range("M86").Value = -24600.625
$ws.Range("N86").Value = -73062.336
$ws.Range("H89").Value = 38021.637
$ws.Range("I89").Value = 25723.625
$ws.Range("J89").Value = 70816.336
$ws.Range("K89").Value = 128618.125
$ws.Range("L89").Value = 354081.68
$ws.Range("M89").Value = -123002.125
$ws.Range("N89").Value = -365313.68
$ws.Range("H92").Value = 251241
$ws.Range("J92").Value = 251241
$ws.Range("L92").Value = 251241
$ws.Range("N92").Value = -256233
$ws.Range("H94").Value = 1996.3889
$ws.Range("I94").Value = 1572.1
$ws.Range("K94").Value = 1572.1
$ws.Range("M94").Value = -1121.1
$ws.Range("H105").Value = 7826.6313
$ws.Range("I105").Value = 13169.7
$ws.Range("J105").Value = 5918.393
$ws.Range("K105").Value = 13169.7
$ws.Range("L105").Value = 5918.393
$ws.Range("M105").Value = -11422.7
$ws.Range("N105").Value = -9412.393
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 1677
$ws.Range("I134").Value = 1631.25
$ws.Range("K134").Value = 4893.75
$ws.Range("M134").Value = -2358.75
$ws.Range("H135").Value = 81924.5
$ws.Range("J135").Value = 81924.5
$ws.Range("L135").Value = 81924.5
$ws.Range("N135").Value = -92064.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4707.4165
$ws.Range("I99").Value = 4680.8184
$ws.Range("K99").Value = 4680.8184
$ws.Range("M99").Value = -3182.8184
$ws.Range("H105").Value = 100011
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H109").Value = 39999.125
$ws.Range("J109").Value = 34999
$ws.Range("L109").Value = 34999
$ws.Range("N109").Value = -37079
$ws.Range("H126").Value = 4707.4165
$ws.Range("I126").Value = 4680.8184
$ws.Range("K126").Value = 14042.4552
$ws.Range("M126").Value = -11572.4552
$ws.Range("H132").Value = 4321.0527
$ws.Range("I132").Value = 4264.8237
$ws.Range("K132").Value = 12794.4711
$ws.Range("M132").Value = -10264.4711
$ws.Range("H134").Value = 3230.7083
$ws.Range("I134").Value = 2775.6843
$ws.Range("K134").Value = 8327.052899999999
$ws.Range("M134").Value = -5792.052899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 200
$ws.Range("I11").Value = 200
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 600
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -460
$ws.Range("N11").ClearContents()
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -1331
$ws.Range("N17").ClearContents()
$ws.Range("H87").Value = 8336.666999999999
$ws.Range("I87").Value = 8336.666999999999
$ws.Range("K87").Value = 25010.001
$ws.Range("M87").Value = -23762.001
$ws.Range("H90").Value = 8336.666999999999
$ws.Range("I90").Value = 8336.666999999999
$ws.Range("K90").Value = 75030.003
$ws.Range("M90").Value = -68790.003
$ws.Range("H93").Value = 1000
$ws.Range("J93").Value = 1000
$ws.Range("L93").Value = 3000
$ws.Range("N93").Value = -6744
$ws.Range("H98").Value = 691.8
$ws.Range("J98").Value = 755.3333
$ws.Range("L98").Value = 2265.9999
$ws.Range("N98").Value = -5261.9999
$ws.Range("H113").Value = 1188.1666
$ws.Range("J113").Value = 3108.5
$ws.Range("L113").Value = 9325.5
$ws.Range("N113").Value = -13665.5
$ws.Range("H121").Value = 1968.5
$ws.Range("J121").Value = 1437.5
$ws.Range("L121").Value = 4312.5
$ws.Range("N121").Value = -6932.5
$ws.Range("H131").Value = 55556856
$ws.Range("I131").Value = 83334430
$ws.Range("J131").Value = 1710
$ws.Range("K131").Value = 250003290
$ws.Range("L131").Value = 5130
$ws.Range("M131").Value = -249998250
$ws.Range("N131").Value = -15210

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7379.857
$ws.Range("I70").Value = 4768
$ws.Range("J70").Value = 8830.888999999999
$ws.Range("K70").Value = 4768
$ws.Range("L70").Value = 8830.888999999999
$ws.Range("M70").Value = -4498
$ws.Range("N70").Value = -9370.888999999999
$ws.Range("H73").Value = 7379.857
$ws.Range("I73").Value = 4768
$ws.Range("J73").Value = 8830.888999999999
$ws.Range("K73").Value = 4768
$ws.Range("L73").Value = 8830.888999999999
$ws.Range("M73").Value = -3832
$ws.Range("N73").Value = -10702.889
$ws.Range("H102").Value = 9199.429
$ws.Range("I102").Value = 1879.2
$ws.Range("J102").Value = 27500
$ws.Range("K102").Value = 1879.2
$ws.Range("L102").Value = 27500
$ws.Range("M102").Value = -257.2
$ws.Range("N102").Value = -30744
$ws.Range("H126").Value = 1888.7
$ws.Range("I126").Value = 1765.2222
$ws.Range("K126").Value = 5295.6666
$ws.Range("M126").Value = -2825.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 833.4737
$ws.Range("I55").Value = 944.5833
$ws.Range("J55").Value = 643
$ws.Range("K55").Value = 944.5833
$ws.Range("L55").Value = 643
$ws.Range("M55").Value = -771.5833
$ws.Range("N55").Value = -989
$ws.Range("H93").Value = 435237.84
$ws.Range("I93").Value = 1615.2593
$ws.Range("J93").Value = 1898714
$ws.Range("K93").Value = 1615.2593
$ws.Range("L93").Value = 1898714
$ws.Range("M93").Value = -367.2592999999999
$ws.Range("N93").Value = -1901210
$ws.Range("H132").Value = 1721.4445
$ws.Range("I132").Value = 1556.8846
$ws.Range("K132").Value = 4670.6538
$ws.Range("M132").Value = -2140.6538
$ws.Range("H136").Value = 7391.1665
$ws.Range("I136").Value = 3736.125
$ws.Range("K136").Value = 11208.375
$ws.Range("M136").Value = -8658.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 12000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H86").Value = 10032908
$ws.Range("J86").Value = 29883
$ws.Range("L86").Value = 29883
$ws.Range("N86").Value = -32129
$ws.Range("H87").Value = 66761.664
$ws.Range("I87").Value = 66761.664
$ws.Range("K87").Value = 66761.664
$ws.Range("M87").Value = -65513.664
$ws.Range("H89").Value = 10032908
$ws.Range("J89").Value = 29883
$ws.Range("L89").Value = 149415
$ws.Range("N89").Value = -160647
$ws.Range("H90").Value = 66761.664
$ws.Range("I90").Value = 66761.664
$ws.Range("K90").Value = 200284.992
$ws.Range("M90").Value = -194044.992
$ws.Range("H96").Value = 6463
$ws.Range("I96").Value = 7466.6665
$ws.Range("J96").Value = 5860.8
$ws.Range("K96").Value = 7466.6665
$ws.Range("L96").Value = 5860.8
$ws.Range("M96").Value = -6093.6665
$ws.Range("N96").Value = -8606.799999999999
$ws.Range("H104").Value = 44950
$ws.Range("J104").Value = 44950
$ws.Range("L104").Value = 44950
$ws.Range("N104").Value = -51938
$ws.Range("H132").Value = 2445.1765
$ws.Range("J132").Value = 4500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560
